$wb = $excel.ActiveWorkbook

# Sheet "space" (1st sheet): change B2 value, update selection, no longer the active tab
$wsSpace = $wb.Worksheets.Item("space")
$wsSpace.Activate()
$wsSpace.Range("B2").Value = 1
$wsSpace.Range("A2:F9").Select()

# Sheet "binned_space" (3rd sheet): change A2 value, update selection, becomes the active tab
$wsBinned = $wb.Worksheets.Item("binned_space")
$wsBinned.Activate()
$wsBinned.Range("A2").Value = 0.75
$wsBinned.Range("F12").Select()
